$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.541.63"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.913.12"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.49"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  +5.82%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.27"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "49.03"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "2.188.16"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.38"
$ws.Range("E14").Value = "  +8.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.703"
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").Value = "1.933.27"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.87"
$ws.Range("E17").Value = "  +3.55%  "
$ws.Range("D18").Value = "35.578.56"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.56"
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("D20").Value = "0.0₃0833"
$ws.Range("E20").Value = "  +4.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "245.03"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.72"
$ws.Range("E22").Value = "  +4.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.85"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +18.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.75"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +6.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.32"
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  +3.90%  "
$ws.Range("E31").Value = "  +25.56%  "
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("E34").Value = "  +5.18%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0207"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "92.72"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0635"
$ws.Range("E42").Value = "  +15.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.61"
$ws.Range("E43").Value = "  +4.99%  "
$ws.Range("D44").Value = "1.352.33"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.95"
$ws.Range("E46").Value = "  +40.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.66"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.79"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.61"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").Value = "2.099.68"
$ws.Range("E51").Value = "  +2.98%  "
